# Apply updated NATMI ligand-receptor statistics (Vcam1-Itgb1) per Dr Hou advice.
# Ligand/Receptor-expressing cell counts change from 1 to 3 for every data row,
# which cascades into recalculated expression, specificity and edge-weight values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.664768333333334
$ws.Range("H2").Value = 19.994305
$ws.Range("I2").Value = 0.06516174319532789
$ws.Range("J2").Value = 0.0651617431953279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 1023.235332970984
$ws.Range("R2").Value = 9209.11799673886
$ws.Range("S2").Value = 0.02067065358645799
$ws.Range("T2").Value = 0.02067065358645799

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.664768333333334
$ws.Range("H3").Value = 19.994305
$ws.Range("I3").Value = 0.06516174319532789
$ws.Range("J3").Value = 0.0651617431953279
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 1125.010913008882
$ws.Range("R3").Value = 10125.09821707994
$ws.Range("S3").Value = 0.02272664959317903
$ws.Range("T3").Value = 0.02272664959317904

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.664768333333334
$ws.Range("H4").Value = 19.994305
$ws.Range("I4").Value = 0.06516174319532789
$ws.Range("J4").Value = 0.0651617431953279
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 453.8062307584278
$ws.Range("R4").Value = 4084.25607682585
$ws.Range("S4").Value = 0.009167462351155613
$ws.Range("T4").Value = 0.009167462351155616

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.664768333333334
$ws.Range("H5").Value = 19.994305
$ws.Range("I5").Value = 0.06516174319532789
$ws.Range("J5").Value = 0.0651617431953279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 623.5735401924218
$ws.Range("R5").Value = 5612.161861731795
$ws.Range("S5").Value = 0.01259697766453526
$ws.Range("T5").Value = 0.01259697766453526

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.25592399999999
$ws.Range("H6").Value = 141.767772
$ws.Range("I6").Value = 0.4620233187619072
$ws.Range("J6").Value = 0.4620233187619072
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 7255.155574898681
$ws.Range("R6").Value = 65296.40017408813
$ws.Range("S6").Value = 0.1465633591533168
$ws.Range("T6").Value = 0.1465633591533168

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 47.25592399999999
$ws.Range("H7").Value = 141.767772
$ws.Range("I7").Value = 0.4620233187619072
$ws.Range("J7").Value = 0.4620233187619072
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 7976.785920438598
$ws.Range("R7").Value = 71791.07328394738
$ws.Range("S7").Value = 0.1611412088517054
$ws.Range("T7").Value = 0.1611412088517054

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.25592399999999
$ws.Range("H8").Value = 141.767772
$ws.Range("I8").Value = 0.4620233187619072
$ws.Range("J8").Value = 0.4620233187619072
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 3217.671144575426
$ws.Range("R8").Value = 28959.04030117884
$ws.Range("S8").Value = 0.06500104466832994
$ws.Range("T8").Value = 0.06500104466832995

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.25592399999999
$ws.Range("H9").Value = 141.767772
$ws.Range("I9").Value = 0.4620233187619072
$ws.Range("J9").Value = 0.4620233187619072
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 4421.390564524851
$ws.Range("R9").Value = 39792.51508072366
$ws.Range("S9").Value = 0.08931770608855505
$ws.Range("T9").Value = 0.08931770608855505

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 20.98736333333333
$ws.Range("H10").Value = 62.96209
$ws.Range("I10").Value = 0.2051944060881897
$ws.Range("J10").Value = 0.2051944060881898
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 3222.169268984298
$ws.Range("R10").Value = 28999.52342085868
$ws.Range("S10").Value = 0.06509191249555264
$ws.Range("T10").Value = 0.06509191249555267

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 20.98736333333333
$ws.Range("H11").Value = 62.96209
$ws.Range("I11").Value = 0.2051944060881897
$ws.Range("J11").Value = 0.2051944060881898
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 3542.66069042397
$ws.Range("R11").Value = 31883.94621381573
$ws.Range("S11").Value = 0.0715662463428562
$ws.Range("T11").Value = 0.07156624634285622

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 20.98736333333333
$ws.Range("H12").Value = 62.96209
$ws.Range("I12").Value = 0.2051944060881897
$ws.Range("J12").Value = 0.2051944060881898
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 1429.036355280811
$ws.Range("R12").Value = 12861.3271975273
$ws.Range("S12").Value = 0.02886834974384312
$ws.Range("T12").Value = 0.02886834974384313

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 20.98736333333333
$ws.Range("H13").Value = 62.96209
$ws.Range("I13").Value = 0.2051944060881897
$ws.Range("J13").Value = 0.2051944060881898
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 1963.633812688857
$ws.Range("R13").Value = 17672.70431419971
$ws.Range("S13").Value = 0.03966789750593775
$ws.Range("T13").Value = 0.03966789750593776

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 27.37233166666667
$ws.Range("H14").Value = 82.116995
$ws.Range("I14").Value = 0.2676205319545753
$ws.Range("J14").Value = 0.2676205319545753
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 4202.447182905416
$ws.Range("R14").Value = 37822.02464614874
$ws.Range("S14").Value = 0.08489477164652148
$ws.Range("T14").Value = 0.0848947716465215

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 27.37233166666667
$ws.Range("H15").Value = 82.116995
$ws.Range("I15").Value = 0.2676205319545753
$ws.Range("J15").Value = 0.2676205319545753
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 4620.441446626719
$ws.Range("R15").Value = 41583.97301964047
$ws.Range("S15").Value = 0.0933387867700245
$ws.Range("T15").Value = 0.0933387867700245

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 27.37233166666667
$ws.Range("H16").Value = 82.116995
$ws.Range("I16").Value = 0.2676205319545753
$ws.Range("J16").Value = 0.2676205319545753
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 1863.790913570572
$ws.Range("R16").Value = 16774.11822213515
$ws.Range("S16").Value = 0.03765094410896171
$ws.Range("T16").Value = 0.03765094410896172

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 27.37233166666667
$ws.Range("H17").Value = 82.116995
$ws.Range("I17").Value = 0.2676205319545753
$ws.Range("J17").Value = 0.2676205319545753
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 2561.028516975879
$ws.Range("R17").Value = 23049.2566527829
$ws.Range("S17").Value = 0.0517360294290676
$ws.Range("T17").Value = 0.0517360294290676

Write-Host "Applied NATMI recalculation updates to rows 2-17"
